$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "fordern"
$ws.Range("B3").Value = "flower/flower013.jpg"
$ws.Range("A4").Value = "lächeln"
$ws.Range("B4").Value = "dog/dog023.jpg"
$ws.Range("C4").Value = "dog"
$ws.Range("A6").Value = "leisten"
$ws.Range("B6").Value = "dog/dog025.jpg"
$ws.Range("C6").Value = "dog"
$ws.Range("A7").Value = "schneiden"
$ws.Range("B7").Value = "flower/flower004.jpg"
$ws.Range("C7").Value = "flower"
$ws.Range("A9").Value = "wandern"
$ws.Range("B9").Value = "dog/dog002.jpg"
$ws.Range("C9").Value = "dog"
$ws.Range("A10").Value = "biegen"
$ws.Range("B10").Value = "flower/flower017.jpg"
$ws.Range("A12").Value = "brauchen"
$ws.Range("B12").Value = "flower/flower011.jpg"
$ws.Range("A13").Value = "passen"
$ws.Range("B13").Value = "flower/flower027.jpg"
$ws.Range("A15").Value = "zögern"
$ws.Range("B15").Value = "flower/flower033.jpg"
$ws.Range("A16").Value = "helfen"
$ws.Range("B16").Value = "flower/flower001.jpg"
$ws.Range("A18").Value = "dringen"
$ws.Range("B18").Value = "flower/flower026.jpg"
$ws.Range("C18").Value = "flower"
$ws.Range("A19").Value = "kümmern"
$ws.Range("B19").Value = "flower/flower030.jpg"
$ws.Range("A21").Value = "zählen"
$ws.Range("B21").Value = "flower/flower012.jpg"
$ws.Range("C21").Value = "flower"
$ws.Range("A22").Value = "landen"
$ws.Range("B22").Value = "dog/dog014.jpg"
$ws.Range("A24").Value = "ruhen"
$ws.Range("B24").Value = "dog/dog015.jpg"
$ws.Range("A25").Value = "faulen"
$ws.Range("B25").Value = "dog/dog019.jpg"
$ws.Range("C25").Value = "dog"
$ws.Range("A27").Value = "grüßen"
$ws.Range("B27").Value = "dog/dog017.jpg"
$ws.Range("A28").Value = "schützen"
$ws.Range("B28").Value = "flower/flower006.jpg"
$ws.Range("C28").Value = "flower"
$ws.Range("A30").Value = "schicken"
$ws.Range("B30").Value = "dog/dog021.jpg"
$ws.Range("A31").Value = "antun"
$ws.Range("B31").Value = "flower/flower023.jpg"
$ws.Range("A33").Value = "treten"
$ws.Range("B33").Value = "flower/flower003.jpg"
$ws.Range("C33").Value = "flower"
$ws.Range("A34").Value = "wundern"
$ws.Range("B34").Value = "dog/dog030.jpg"
$ws.Range("C34").Value = "dog"
$ws.Range("A36").Value = "fühlen"
$ws.Range("B36").Value = "dog/dog001.jpg"
$ws.Range("A37").Value = "lassen"
$ws.Range("B37").Value = "dog/dog010.jpg"
$ws.Range("C37").Value = "dog"
$ws.Range("A39").Value = "wehtun"
$ws.Range("B39").Value = "flower/flower024.jpg"
$ws.Range("C39").Value = "flower"
$ws.Range("A40").Value = "segnen"
$ws.Range("B40").Value = "dog/dog013.jpg"
$ws.Range("C40").Value = "dog"
$ws.Range("A42").Value = "rasen"
$ws.Range("B42").Value = "dog/dog003.jpg"
$ws.Range("A43").Value = "fügen"
$ws.Range("B43").Value = "dog/dog005.jpg"
$ws.Range("A45").Value = "zeugen"
$ws.Range("B45").Value = "flower/flower002.jpg"
$ws.Range("A46").Value = "posten"
$ws.Range("B46").Value = "dog/dog028.jpg"
$ws.Range("A48").Value = "leugnen"
$ws.Range("B48").Value = "flower/flower031.jpg"
$ws.Range("C48").Value = "flower"
$ws.Range("A49").Value = "spielen"
$ws.Range("B49").Value = "dog/dog009.jpg"
